$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data refresh swapped which row holds ASML vs. Taiwan Semiconductor
# (TSM), so row 2 now carries ASML's name/ticker/metrics and row 3 carries
# TSM's - plus every row's metrics (D:K, N) were refreshed with the latest
# pull.

# Row 2 -> ASML Holding N.V. - New York Re / ASML
$ws.Range("B2").Value = "ASML Holding N.V. - New York Re"
$ws.Range("C2").Value = "ASML"
$ws.Range("D2").Value = 1110.08
$ws.Range("E2").Value = 62.7
$ws.Range("F2").Value = 6.64
$ws.Range("H2").Value = 53
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 59.7
$ws.Range("N2").Value = 54.85170003294819

# Row 3 -> Taiwan Semiconductor Manufactur / TSM
$ws.Range("B3").Value = "Taiwan Semiconductor Manufactur"
$ws.Range("C3").Value = "TSM"
$ws.Range("D3").Value = 292.93
$ws.Range("E3").Value = 60.3
$ws.Range("F3").Value = 1.02
$ws.Range("H3").Value = 70
$ws.Range("J3").Value = 83
$ws.Range("K3").Value = 58.5
$ws.Range("N3").Value = 54.85170003294819

# Row 4 -> QUALCOMM Incorporated / QCOM refreshed metrics
$ws.Range("D4").Value = 174.35
$ws.Range("E4").Value = 49.8
$ws.Range("F4").Value = 5.58
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 40
$ws.Range("K4").Value = 50.5
$ws.Range("N4").Value = 54.85170003294819

# Row 5 -> NVIDIA Corporation / NVDA refreshed metrics
$ws.Range("D5").Value = 183.38
$ws.Range("E5").Value = 46.5
$ws.Range("F5").Value = 1.73
$ws.Range("H5").Value = 73
$ws.Range("K5").Value = 49.5
$ws.Range("N5").Value = 54.85170003294819

# Row 6 -> Advanced Micro Devices, Inc. / AMD refreshed metrics
$ws.Range("D6").Value = 215.98
$ws.Range("F6").Value = 0.8100000000000001
$ws.Range("K6").Value = 46.5
$ws.Range("N6").Value = 54.85170003294819
